{"js": "// Office.js (Word JavaScript API) edit script.\n// The document body consists of 101 paragraphs in document order:\n//   paragraph 0  -> the title line with the date\n//   paragraphs 1..100 -> one arithmetic-equation paragraph per table cell\n// The commit replaces each paragraph's text with a new value, preserving\n// the existing run formatting (font/size) by replacing the text in place.\nconst newTexts = [\n  \"2025-02-21 Friday\",\n  \"60+13=73\",\n  \"7+65=72\",\n  \"16+77=93\",\n  \"96-90=6\",\n  \"35+29=64\",\n  \"2+16=18\",\n  \"27+15=42\",\n  \"45-17=28\",\n  \"44-38=6\",\n  \"61-56=5\",\n  \"5+39=44\",\n  \"22+57=79\",\n  \"37+55=92\",\n  \"65+0=65\",\n  \"75-72=3\",\n  \"10+51=61\",\n  \"14+38=52\",\n  \"41+6=47\",\n  \"82+7=89\",\n  \"79-0=79\",\n  \"45+12=57\",\n  \"54+35=89\",\n  \"32+2=34\",\n  \"77-53=24\",\n  \"58+27=85\",\n  \"21+4=25\",\n  \"49+10=59\",\n  \"17+16=33\",\n  \"67-58=9\",\n  \"46+30=76\",\n  \"70-33=37\",\n  \"3+76=79\",\n  \"64+32=96\",\n  \"72-20=52\",\n  \"92-50=42\",\n  \"2+79=81\",\n  \"5+2=7\",\n  \"22-6=16\",\n  \"47-22=25\",\n  \"88-86=2\",\n  \"22+68=90\",\n  \"13+26=39\",\n  \"95-60=35\",\n  \"7+70=77\",\n  \"93-88=5\",\n  \"4+77=81\",\n  \"94-36=58\",\n  \"52+38=90\",\n  \"75-29=46\",\n  \"93-28=65\",\n  \"29+5=34\",\n  \"18+70=88\",\n  \"52-48=4\",\n  \"12+19=31\",\n  \"9+35=44\",\n  \"44-3=41\",\n  \"95-7=88\",\n  \"65-61=4\",\n  \"85+2=87\",\n  \"98-98=0\",\n  \"89-41=48\",\n  \"30+2=32\",\n  \"19+37=56\",\n  \"99-44=55\",\n  \"4+94=98\",\n  \"3+58=61\",\n  \"10+55=65\",\n  \"77-64=13\",\n  \"94-20=74\",\n  \"45+30=75\",\n  \"16+7=23\",\n  \"99-38=61\",\n  \"96-41=55\",\n  \"40+13=53\",\n  \"55+13=68\",\n  \"32+51=83\",\n  \"73-59=14\",\n  \"88-30=58\",\n  \"26+13=39\",\n  \"51-18=33\",\n  \"36+40=76\",\n  \"11+1=12\",\n  \"36+62=98\",\n  \"23-18=5\",\n  \"98-93=5\",\n  \"14+19=33\",\n  \"11+8=19\",\n  \"10+3=13\",\n  \"96-52=44\",\n  \"20-13=7\",\n  \"55-0=55\",\n  \"22-15=7\",\n  \"51-29=22\",\n  \"31+30=61\",\n  \"63-4=59\",\n  \"23+52=75\",\n  \"28-7=21\",\n  \"15-14=1\",\n  \"42-39=3\",\n  \"92-85=7\"\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nif (paragraphs.items.length !== newTexts.length) {\n  throw new Error(\n    \"Unexpected paragraph count: expected \" + newTexts.length +\n    \" got \" + paragraphs.items.length\n  );\n}\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].insertText(newTexts[i], Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Title paragraph holds the date; the 20x5 table holds one arithmetic\n# equation per cell. Replace the title text and every cell's text, in\n# document/row-major order, preserving existing run formatting by\n# writing straight into each Range.Text.\n\n$d = $word.ActiveDocument\n\n# 1) Title paragraph with the date.\n$d.Paragraphs.Item(1).Range.Text = \"2025-02-21 Friday\"\n\n# 2) Table cells, row-major (row 1 col 1..5, row 2 col 1..5, ...).\n$newCellTexts = @(\n  \"60+13=73\",\n  \"7+65=72\",\n  \"16+77=93\",\n  \"96-90=6\",\n  \"35+29=64\",\n  \"2+16=18\",\n  \"27+15=42\",\n  \"45-17=28\",\n  \"44-38=6\",\n  \"61-56=5\",\n  \"5+39=44\",\n  \"22+57=79\",\n  \"37+55=92\",\n  \"65+0=65\",\n  \"75-72=3\",\n  \"10+51=61\",\n  \"14+38=52\",\n  \"41+6=47\",\n  \"82+7=89\",\n  \"79-0=79\",\n  \"45+12=57\",\n  \"54+35=89\",\n  \"32+2=34\",\n  \"77-53=24\",\n  \"58+27=85\",\n  \"21+4=25\",\n  \"49+10=59\",\n  \"17+16=33\",\n  \"67-58=9\",\n  \"46+30=76\",\n  \"70-33=37\",\n  \"3+76=79\",\n  \"64+32=96\",\n  \"72-20=52\",\n  \"92-50=42\",\n  \"2+79=81\",\n  \"5+2=7\",\n  \"22-6=16\",\n  \"47-22=25\",\n  \"88-86=2\",\n  \"22+68=90\",\n  \"13+26=39\",\n  \"95-60=35\",\n  \"7+70=77\",\n  \"93-88=5\",\n  \"4+77=81\",\n  \"94-36=58\",\n  \"52+38=90\",\n  \"75-29=46\",\n  \"93-28=65\",\n  \"29+5=34\",\n  \"18+70=88\",\n  \"52-48=4\",\n  \"12+19=31\",\n  \"9+35=44\",\n  \"44-3=41\",\n  \"95-7=88\",\n  \"65-61=4\",\n  \"85+2=87\",\n  \"98-98=0\",\n  \"89-41=48\",\n  \"30+2=32\",\n  \"19+37=56\",\n  \"99-44=55\",\n  \"4+94=98\",\n  \"3+58=61\",\n  \"10+55=65\",\n  \"77-64=13\",\n  \"94-20=74\",\n  \"45+30=75\",\n  \"16+7=23\",\n  \"99-38=61\",\n  \"96-41=55\",\n  \"40+13=53\",\n  \"55+13=68\",\n  \"32+51=83\",\n  \"73-59=14\",\n  \"88-30=58\",\n  \"26+13=39\",\n  \"51-18=33\",\n  \"36+40=76\",\n  \"11+1=12\",\n  \"36+62=98\",\n  \"23-18=5\",\n  \"98-93=5\",\n  \"14+19=33\",\n  \"11+8=19\",\n  \"10+3=13\",\n  \"96-52=44\",\n  \"20-13=7\",\n  \"55-0=55\",\n  \"22-15=7\",\n  \"51-29=22\",\n  \"31+30=61\",\n  \"63-4=59\",\n  \"23+52=75\",\n  \"28-7=21\",\n  \"15-14=1\",\n  \"42-39=3\",\n  \"92-85=7\"\n)\n\n$table = $d.Tables.Item(1)\n$rows = $table.Rows.Count\n$cols = $table.Columns.Count\n\nif (($rows * $cols) -ne $newCellTexts.Count) {\n  throw \"Unexpected table size: expected $($newCellTexts.Count) cells, got $($rows * $cols)\"\n}\n\n$i = 0\nfor ($r = 1; $r -le $rows; $r++) {\n  for ($c = 1; $c -le $cols; $c++) {\n    $table.Cell($r, $c).Range.Text = $newCellTexts[$i]\n    $i++\n  }\n}\n"}
